# Update countries & provincias Spain
# Refreshes the COVID dashboard data in the "Pais" worksheet:
#  - bumps the "last updated" timestamp
#  - refreshes case/death/etc. counters for the countries whose figures moved
#  - re-labels four rows whose country swapped rank with its neighbour
#    (Uzbekistan/Libano, Angola/Cabo Verde, Congo/Siria, Belice/Republica de Chipre)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 21:12"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8487233
$ws.Range("C4").Value = 30580
$ws.Range("D4").Value = 5521771
$ws.Range("E4").Value = 2739717
$ws.Range("G4").Value = 523
$ws.Range("H4").Value = 225745

# Row 5 - India
$ws.Range("B5").Value = 7648373
$ws.Range("C5").Value = 53637
$ws.Range("D5").Value = 6791113
$ws.Range("E5").Value = 741321
$ws.Range("G5").Value = 703
$ws.Range("H5").Value = 115939

# Row 11 - Francia
$ws.Range("B11").Value = 930745
$ws.Range("C11").Value = 20468
$ws.Range("D11").Value = 106839
$ws.Range("E11").Value = 790021
$ws.Range("G11").Value = 262
$ws.Range("H11").Value = 33885

# Row 21 - Alemania
$ws.Range("B21").Value = 380022
$ws.Range("C21").Value = 6291
$ws.Range("E21").Value = 71773
$ws.Range("G21").Value = 50
$ws.Range("H21").Value = 9949

# Row 31 - Canada
$ws.Range("B31").Value = 203155
$ws.Range("C31").Value = 1718
$ws.Range("D31").Value = 171354
$ws.Range("E31").Value = 22009

# Row 62 - now Libano (was Uzbekistan)
$ws.Range("A62").Value = "Libano"
$ws.Range("B62").Value = 64336
$ws.Range("C62").Value = 1392
$ws.Range("D62").Value = 29498
$ws.Range("E62").Value = 34307
$ws.Range("G62").Value = 5
$ws.Range("H62").Value = 531

# Row 63 - now Uzbekistan (was Libano)
$ws.Range("A63").Value = "Uzbekistan"
$ws.Range("B63").Value = 63831
$ws.Range("C63").Value = 308
$ws.Range("D63").Value = 60910
$ws.Range("E63").Value = 2387
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 534

# Row 118 - now Cabo Verde (was Angola)
$ws.Range("A118").Value = "Cabo Verde"
$ws.Range("B118").Value = 7901
$ws.Range("C118").Value = 101
$ws.Range("D118").Value = 6792
$ws.Range("E118").Value = 1022
$ws.Range("H118").Value = 87

# Row 119 - now Angola (was Cabo Verde)
$ws.Range("A119").Value = "Angola"
$ws.Range("B119").Value = 7829
$ws.Range("D119").Value = 3031
$ws.Range("E119").Value = 4550
$ws.Range("H119").Value = 248

# Row 132 - now Siria (was Congo)
$ws.Range("A132").Value = "Siria"
$ws.Range("B132").Value = 5180
$ws.Range("C132").Value = 46
$ws.Range("D132").Value = 1596
$ws.Range("E132").Value = 3330
$ws.Range("G132").Value = 3
$ws.Range("H132").Value = 254

# Row 133 - now Congo (was Siria)
$ws.Range("A133").Value = "Congo"
$ws.Range("B133").Value = 5156
$ws.Range("D133").Value = 3887
$ws.Range("E133").Value = 1177
$ws.Range("H133").Value = 92

# Row 141 - Aruba
$ws.Range("B141").Value = 4355
$ws.Range("C141").Value = 21
$ws.Range("D141").Value = 4065
$ws.Range("E141").Value = 256

# Row 153 - now Republica de Chipre (was Belice)
$ws.Range("A153").Value = "Republica de Chipre"
$ws.Range("B153").Value = 2839
$ws.Range("C153").Value = 152
$ws.Range("D153").Value = 1444
$ws.Range("E153").Value = 1370
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 25

# Row 154 - now Belice (was Republica de Chipre)
$ws.Range("A154").Value = "Belice"
$ws.Range("B154").Value = 2833
$ws.Range("C154").Value = 20
$ws.Range("D154").Value = 1692
$ws.Range("E154").Value = 1096
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 45

# Row 162 - Yemen
$ws.Range("B162").Value = 2057
$ws.Range("C162").Value = 1
$ws.Range("E162").Value = 122

# Row 189 - Monaco
$ws.Range("B189").Value = 271
$ws.Range("C189").Value = 3
$ws.Range("D189").Value = 227
$ws.Range("E189").Value = 42

# Row 200 - San Vicente y las Granadinas
$ws.Range("B200").Value = 68
$ws.Range("C200").Value = 1
$ws.Range("E200").Value = 4
